$d = $word.ActiveDocument

# 1. Remove the stray "_GoBack" bookmark (start+end) near the top of the
#    document. It is a hidden bookmark so it isn't enumerated by the
#    normal Bookmarks collection, but it can still be reached directly
#    by name.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# 2. Fix the typing table for the "a"/"አ" row: the romanization under
#    አ should read "ee" (was "a"), and the romanization under ኣ should
#    read "a" (was "aa"). Each Find is scoped to its own table cell so
#    the edits can't clobber each other or match unrelated cells.
$t = $d.Tables.Item(1)

$cellA = $t.Cell(16, 2)
$rngA = $cellA.Range
[void]$rngA.MoveEnd(1, -1)
[void]$rngA.Find.Execute("a", $true, $false, $false, $false, $false, $true, 1, $false, "ee", 1)

$cellAa = $t.Cell(16, 5)
$rngAa = $cellAa.Range
[void]$rngAa.MoveEnd(1, -1)
[void]$rngAa.Find.Execute("aa", $true, $false, $false, $false, $false, $true, 1, $false, "a", 1)
